$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.341.38'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.878.80'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7107'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.55'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08045'
$ws.Range('E8').Value = '  +3.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3171'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.05'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08321'
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.890.36'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.260'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.72'
$ws.Range('E14').Value = '  +3.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7174'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.400'
$ws.Range('E16').Value = '  +5.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008678'
$ws.Range('E17').Value = '  +4.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.351.94'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.60'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.34'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.136.36'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.819'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1570'
$ws.Range('E25').Value = '  -1.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.100'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.10'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.57'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.431'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.340'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.197'
$ws.Range('E32').Value = '  -7.42%  '
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.944'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7729'
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.188'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.685'
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.265.80'
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.751'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.495'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '113.37'
$ws.Range('E42').Value = '  +2.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '74.65'
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9083'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('E45').Value = '  +7.30%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.028.54'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.811'
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5223'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.509'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4374'
$ws.Range('E51').Value = '  +1.21%  '
